$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.4165321785012566
$ws.Range("D2").Value = 0.1770283408321944
$ws.Range("E2").Value = 0.1759721268987846
$ws.Range("F2").Value = 3.387674600170385
$ws.Range("G2").Value = 3.717185640561581
$ws.Range("H2").Value = 2.478075300866749
$ws.Range("L2").Value = 0.1535809029467927
$ws.Range("C3").Value = 0.408617831278633
$ws.Range("D3").Value = 0.1696488696270393
$ws.Range("E3").Value = 0.172217812073022
$ws.Range("F3").Value = 3.20830773157536
$ws.Range("G3").Value = 3.493995001652593
$ws.Range("H3").Value = 2.383840216937926
$ws.Range("L3").Value = 0.1499384042475285
$ws.Range("C4").Value = 0.4040558284525844
$ws.Range("D4").Value = 0.1652659502851179
$ws.Range("E4").Value = 0.1700421457223378
$ws.Range("F4").Value = 3.100954056404646
$ws.Range("G4").Value = 3.360036834047548
$ws.Range("H4").Value = 2.327943079853128
$ws.Range("L4").Value = 0.1478181069541549
$ws.Range("C5").Value = 0.402270600351784
$ws.Range("D5").Value = 0.1635162936335064
$ws.Range("E5").Value = 0.1691876484232715
$ws.Range("F5").Value = 3.05788877680348
$ws.Range("G5").Value = 3.306202523923787
$ws.Range("H5").Value = 2.305649011294008
$ws.Range("L5").Value = 0.146982845631328
$ws.Range("C6").Value = 0.4019785956473072
$ws.Range("D6").Value = 0.1632279403856529
$ws.Range("E6").Value = 0.1690476856441698
$ws.Range("F6").Value = 3.050778551949605
$ws.Range("G6").Value = 3.297308380218055
$ws.Range("H6").Value = 2.301976091865527
$ws.Range("L6").Value = 0.146845876393435
$ws.Range("C7").Value = 0.4040314545845263
$ws.Range("D7").Value = 0.1652422074195385
$ws.Range("E7").Value = 0.1700304922802758
$ws.Range("F7").Value = 3.100370522109301
$ws.Range("G7").Value = 3.359307776180344
$ws.Range("H7").Value = 2.32764046474864
$ws.Range("L7").Value = 0.1478067263899376
$ws.Range("C8").Value = 0.4137409840663224
$ws.Range("D8").Value = 0.1744526551100023
$ws.Range("E8").Value = 0.1746504803007056
$ws.Range("F8").Value = 3.325242206422672
$ws.Range("G8").Value = 3.639577500156065
$ws.Range("H8").Value = 2.4451698046758
$ws.Range("L8").Value = 0.1523005766433556
$ws.Range("C9").Value = 0.4351876940660304
$ws.Range("D9").Value = 0.1937285989380229
$ws.Range("E9").Value = 0.1847596816637846
$ws.Range("F9").Value = 3.789042000569367
$ws.Range("G9").Value = 4.21460069015302
$ws.Range("H9").Value = 2.691668150337307
$ws.Range("L9").Value = 0.1620563131877475
$ws.Range("C10").Value = 0.4524796716726485
$ws.Range("D10").Value = 0.2086880545117822
$ws.Range("E10").Value = 0.1928588232532178
$ws.Range("F10").Value = 4.144858485178361
$ws.Range("G10").Value = 4.653973240002358
$ws.Range("H10").Value = 2.883189741785941
$ws.Range("L10").Value = 0.1698300199496714
$ws.Range("C11").Value = 0.4606950994637202
$ws.Range("D11").Value = 0.2156796332964177
$ws.Range("E11").Value = 0.1966965571949544
$ws.Range("F11").Value = 4.310259620081922
$ws.Range("G11").Value = 4.857845151148069
$ws.Range("H11").Value = 2.972728262915098
$ws.Range("L11").Value = 0.1735051820266591
$ws.Range("C12").Value = 0.4638575714740796
$ws.Range("D12").Value = 0.2183550493000723
$ws.Range("E12").Value = 0.1981724717683875
$ws.Range("F12").Value = 4.373422961433505
$ws.Range("G12").Value = 4.935647507224985
$ws.Range("H12").Value = 3.006993552294546
$ws.Range("L12").Value = 0.1749174252540087
$ws.Range("C13").Value = 0.4631741670932286
$ws.Range("D13").Value = 0.2177775946002498
$ws.Range("E13").Value = 0.1978535902517109
$ws.Range("F13").Value = 4.359795717886527
$ws.Range("G13").Value = 4.918864256875111
$ws.Range("H13").Value = 2.999597750250075
$ws.Range("L13").Value = 0.1746123506163428
$ws.Range("C14").Value = 0.4609542391178536
$ws.Range("D14").Value = 0.2158991769627789
$ws.Range("E14").Value = 0.1968175242162218
$ws.Range("F14").Value = 4.315445372688544
$ws.Range("G14").Value = 4.864233814428871
$ws.Range("H14").Value = 2.975540024997542
$ws.Range("L14").Value = 0.1736209531485713
$ws.Range("C15").Value = 0.4596012096658626
$ws.Range("D15").Value = 0.2147522513859315
$ws.Range("E15").Value = 0.1961858711810223
$ws.Range("F15").Value = 4.288349091472469
$ws.Range("G15").Value = 4.830850061484455
$ws.Range("H15").Value = 2.960851082899808
$ws.Range("L15").Value = 0.1730163856635585
$ws.Range("C16").Value = 0.4519499071875259
$ws.Range("D16").Value = 0.2082349787647502
$ws.Range("E16").Value = 0.1926111549514857
$ws.Range("F16").Value = 4.134122087581659
$ws.Range("G16").Value = 4.6407323570362
$ws.Range("H16").Value = 2.877387772502686
$ws.Range("L16").Value = 0.1695926820743381
$ws.Range("C17").Value = 0.4473463998458271
$ws.Range("D17").Value = 0.2042853293908991
$ws.Range("E17").Value = 0.1904578833457293
$ws.Range("F17").Value = 4.040429598389409
$ws.Range("G17").Value = 4.525143391300105
$ws.Range("H17").Value = 2.826812389448946
$ws.Range("L17").Value = 0.1675283167677435
$ws.Range("C18").Value = 0.4447313769314292
$ws.Range("D18").Value = 0.2020310800899949
$ws.Range("E18").Value = 0.1892337800419313
$ws.Range("F18").Value = 3.986871930230905
$ws.Range("G18").Value = 4.459034521134242
$ws.Range("H18").Value = 2.797949289413168
$ws.Range("L18").Value = 0.1663539836133907
$ws.Range("C19").Value = 0.4438515699405912
$ws.Range("D19").Value = 0.2012708017278442
$ws.Range("E19").Value = 0.188821775545513
$ws.Range("F19").Value = 3.968794587278524
$ws.Range("G19").Value = 4.43671481326902
$ws.Range("H19").Value = 2.78821529208733
$ws.Range("L19").Value = 0.1659585963332972
$ws.Range("C20").Value = 0.4478330482260162
$ws.Range("D20").Value = 0.2047039592654301
$ws.Range("E20").Value = 0.190685608069515
$ws.Range("F20").Value = 4.050368838828405
$ws.Range("G20").Value = 4.537409057556147
$ws.Range("H20").Value = 2.832172696867929
$ws.Range("L20").Value = 0.1677467188201547
$ws.Range("C21").Value = 0.4616048787614773
$ws.Range("D21").Value = 0.2164501492104591
$ws.Range("E21").Value = 0.1971212224547401
$ws.Range("F21").Value = 4.32845760456712
$ws.Range("G21").Value = 4.880263595686813
$ws.Range("H21").Value = 2.98259652285725
$ws.Range("L21").Value = 0.1739115885415998
$ws.Range("C22").Value = 0.4709061736978128
$ws.Range("D22").Value = 0.2242898467257533
$ws.Range("E22").Value = 0.2014595680059585
$ws.Range("F22").Value = 4.513301266596727
$ws.Range("G22").Value = 5.107851384150877
$ws.Range("H22").Value = 3.083005467025487
$ws.Range("L22").Value = 0.1780606847378436
$ws.Range("C23").Value = 0.4659139666193539
$ws.Range("D23").Value = 0.2200903934246128
$ws.Range("E23").Value = 0.1991318052415565
$ws.Range("F23").Value = 4.414356361456669
$ws.Range("G23").Value = 4.986053475071969
$ws.Range("H23").Value = 3.029219319744527
$ws.Range("L23").Value = 0.1758350602495398
$ws.Range("C24").Value = 0.4476129361772792
$ws.Range("D24").Value = 0.2045146455190263
$ws.Range("E24").Value = 0.1905826106539266
$ws.Range("F24").Value = 4.045874352325768
$ws.Range("G24").Value = 4.531862677884931
$ws.Range("H24").Value = 2.829748638198907
$ws.Range("L24").Value = 0.1676479403671749
$ws.Range("C25").Value = 0.4291212773185293
$ws.Range("D25").Value = 0.18837815320272
$ws.Range("E25").Value = 0.1819092697119373
$ws.Range("F25").Value = 3.661011727456582
$ws.Range("G25").Value = 4.056177817182743
$ws.Range("H25").Value = 2.623202387293077
$ws.Range("L25").Value = 0.1593129888774882
